$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two data rows are removed from the table entirely: "RM 232" (originally
# row 26) and "SC 92" (originally row 28, which becomes row 27 once the
# first row has been removed and everything below it has shifted up).
# Removing both rows shrinks the used range from A1:F35 down to A1:F33 and
# shifts every row below them up by two.
$ws.Range("A26:F26").EntireRow.Delete()
$ws.Range("A27:F27").EntireRow.Delete()

# After the rows above are gone, a handful of values in column D ("C")
# are corrected for the remaining records: some previously-missing values
# are now filled in, and some previously-present values are now cleared.
$ws.Range("D19").Value = -15.5   # RM 125: was blank
$ws.Range("D21").ClearContents() # RM 135: was -14.3
$ws.Range("D23").Value = -13.9   # RM 140: was blank
$ws.Range("D27").ClearContents() # SC 101: was -14.6
$ws.Range("D33").Value = -14.1   # SC 232: was blank
